# Remove the hero-section scratch content (and the two blank paragraphs
# that preceded it) that was appended after the "Our thoughtfully
# prepared party bags..." paragraph, leaving that paragraph followed
# directly by the single trailing blank paragraph.

$d = $word.ActiveDocument

$target = "Our thoughtfully prepared party bags"

$anchor = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "$target*") {
        $anchor = $p
        break
    }
}

$afterScript = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*bootstrap.bundle.min.js*") {
        $afterScript = $p
        break
    }
}

$r = $d.Range($anchor.Range.End, $afterScript.Range.End)
$r.Delete()
